$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh (GitHub Actions scheduled update).
# Columns D (Price) and E (Volume 1h) hold text-formatted values (some
# prices use "." as a thousands separator, e.g. "49.871.39", which Excel
# cannot parse as a number) -- force NumberFormat "@" first on any cell
# whose new value *would* otherwise be auto-parsed as numeric, so the
# stored cell stays text exactly like the source data.

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "49.871.39"
$ws.Range("E2").Value = "  +3.69%  "
$ws.Range("D3").Value = "2.647.44"
$ws.Range("E3").Value = "  +5.84%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("B5").Value = "Solana"
$ws.Range("C5").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D5").Value = "111.35"
$ws.Range("E5").Value = "  +3.71%  "
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "327.24"
$ws.Range("E6").Value = "  +2.09%  "
$ws.Range("E7").Value = "  +1.08%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "0.558"
$ws.Range("E9").Value = "  +3.22%  "
$ws.Range("D10").Value = "40.90"
$ws.Range("E10").Value = "  +3.02%  "
$ws.Range("D11").Value = "20.46"
$ws.Range("E11").Value = "  +1.67%  "
$ws.Range("D12").Value = "0.0821"
$ws.Range("E12").Value = "  +1.22%  "
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("D14").Value = "7.33"
$ws.Range("E14").Value = "  +3.30%  "
$ws.Range("D15").Value = "3.066.50"
$ws.Range("E15").Value = "  +6.01%  "
$ws.Range("D16").Value = "2.650.99"
$ws.Range("E16").Value = "  +5.77%  "
$ws.Range("D17").Value = "0.880"
$ws.Range("E17").Value = "  +5.55%  "
$ws.Range("D18").Value = "49.846.40"
$ws.Range("E18").Value = "  +3.98%  "
$ws.Range("D19").Value = "13.29"
$ws.Range("E19").Value = "  +2.63%  "
$ws.Range("D20").Value = "6.82"
$ws.Range("E20").Value = "  +1.99%  "
$ws.Range("D21").Value = "2.93"
$ws.Range("E21").Value = "  +6.13%  "
$ws.Range("E22").Value = "  +2.11%  "
$ws.Range("D23").Value = "73.01"
$ws.Range("E23").Value = "  +2.11%  "
$ws.Range("D24").Value = "280.39"
$ws.Range("E24").Value = "  +1.63%  "
$ws.Range("D25").Value = "2.59"
$ws.Range("E25").Value = "  +2.45%  "
$ws.Range("D26").Value = "26.96"
$ws.Range("E26").Value = "  +4.28%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("E28").Value = "  -0.90%  "
$ws.Range("D29").Value = "36.76"
$ws.Range("E29").Value = "  +4.73%  "
$ws.Range("D30").Value = "9.96"
$ws.Range("E30").Value = "  +2.68%  "
$ws.Range("E31").Value = "  +2.09%  "
$ws.Range("D32").Value = "49.75"
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("D33").Value = "19.56"
$ws.Range("E33").Value = "  +0.54%  "
$ws.Range("E34").Value = "  +2.61%  "
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("D36").Value = "0.0797"
$ws.Range("E36").Value = "  +1.66%  "
$ws.Range("D37").Value = "2.06"
$ws.Range("E37").Value = "  +6.32%  "
$ws.Range("D38").Value = "4.76"
$ws.Range("E38").Value = "  +3.04%  "
$ws.Range("D39").Value = "3.11"
$ws.Range("E39").Value = "  +8.15%  "
$ws.Range("D40").Value = "126.84"
$ws.Range("E40").Value = "  +4.49%  "
$ws.Range("E41").Value = "  +1.75%  "
$ws.Range("D42").Value = "2.25"
$ws.Range("E42").Value = "  +1.68%  "
$ws.Range("D43").Value = "22.28"
$ws.Range("E43").Value = "  +4.22%  "
$ws.Range("E44").Value = "  +3.72%  "
$ws.Range("E45").Value = "  +7.94%  "
$ws.Range("D46").Value = "2.068.82"
$ws.Range("E46").Value = "  +2.47%  "
$ws.Range("E47").Value = "  +14.20%  "
$ws.Range("E48").Value = "  +7.79%  "
$ws.Range("D49").Value = "9.08"
$ws.Range("E49").Value = "  +1.04%  "
$ws.Range("D50").Value = "5.41"
$ws.Range("E50").Value = "  +4.66%  "
$ws.Range("D51").Value = "81.86"
$ws.Range("E51").Value = "  +1.87%  "
